$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings such as "0.9998" or "13.68" that look like
# plain numbers. Force the range to Text format first so the COM layer
# keeps the new values as text (matching the original inline-string cells)
# instead of silently converting them to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '24.813.52'
$ws.Range("E2").Value = '  +2.35%  '
$ws.Range("D3").Value = '1.703.17'
$ws.Range("E3").Value = '  +1.81%  '
$ws.Range("D4").Value = '0.9998'
$ws.Range("E4").Value = '  -0.30%  '
$ws.Range("D5").Value = '309.23'
$ws.Range("E5").Value = '  +0.64%  '
$ws.Range("D6").Value = '0.9950'
$ws.Range("E6").Value = '  -0.28%  '
$ws.Range("D7").Value = '0.3725'
$ws.Range("E7").Value = '  +0.52%  '
$ws.Range("D8").Value = '49.30'
$ws.Range("E8").Value = '  +4.72%  '
$ws.Range("D9").Value = '0.3436'
$ws.Range("E9").Value = '  +0.01%  '
$ws.Range("D10").Value = '1.198'
$ws.Range("E10").Value = '  +0.79%  '
$ws.Range("D11").Value = '0.07484'
$ws.Range("E11").Value = '  +3.08%  '
$ws.Range("D12").Value = '0.9967'
$ws.Range("E12").Value = '  -0.30%  '
$ws.Range("D13").Value = '20.95'
$ws.Range("E13").Value = '  +3.00%  '
$ws.Range("D14").Value = '6.246'
$ws.Range("E14").Value = '  +2.52%  '
$ws.Range("D15").Value = '6.994'
$ws.Range("E15").Value = '  +3.91%  '
$ws.Range("D16").Value = '1.704.59'
$ws.Range("E16").Value = '  +1.74%  '
$ws.Range("D17").Value = '0.00001128'
$ws.Range("E17").Value = '  +1.78%  '
$ws.Range("D18").Value = '0.06736'
$ws.Range("E18").Value = '  +0.53%  '
$ws.Range("D19").Value = '0.9950'
$ws.Range("E19").Value = '  -0.29%  '
$ws.Range("D20").Value = '84.15'
$ws.Range("E20").Value = '  +3.91%  '
$ws.Range("D21").Value = '17.20'
$ws.Range("E21").Value = '  +4.58%  '
$ws.Range("D22").Value = '6.350'
$ws.Range("E22").Value = '  +4.15%  '
$ws.Range("D23").Value = '13.12'
$ws.Range("E23").Value = '  +9.75%  '
$ws.Range("D24").Value = '24.806.79'
$ws.Range("E24").Value = '  +2.51%  '
$ws.Range("E25").Value = '  +0.63%  '
$ws.Range("D26").Value = '2.767'
$ws.Range("E26").Value = '  +4.21%  '
$ws.Range("D27").Value = '20.33'
$ws.Range("E27").Value = '  +4.19%  '
$ws.Range("D28").Value = '150.37'
$ws.Range("E28").Value = '  -0.85%  '
$ws.Range("D29").Value = '131.67'
$ws.Range("E29").Value = '  +3.44%  '
$ws.Range("D30").Value = '1.892.43'
$ws.Range("E30").Value = '  +1.76%  '
$ws.Range("D31").Value = '1.200'
$ws.Range("E31").Value = '  +23.43%  '
$ws.Range("D32").Value = '6.796'
$ws.Range("E32").Value = '  +7.52%  '
$ws.Range("D33").Value = '4.192'
$ws.Range("E33").Value = '  +4.17%  '
$ws.Range("D34").Value = '1.797'
$ws.Range("E34").Value = '  +2.40%  '
$ws.Range("B35").Value = 'Aptos'
$ws.Range("C35").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D35").Value = '13.68'
$ws.Range("E35").Value = '  +11.45%  '
$ws.Range("B36").Value = 'Stellar'
$ws.Range("C36").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D36").Value = '0.08779'
$ws.Range("E36").Value = '  +4.00%  '
$ws.Range("D37").Value = '5.542'
$ws.Range("E37").Value = '  +3.79%  '
$ws.Range("D38").Value = '0.06628'
$ws.Range("E38").Value = '  +3.73%  '
$ws.Range("D39").Value = '0.02403'
$ws.Range("E39").Value = '  +2.70%  '
$ws.Range("D40").Value = '9.014'
$ws.Range("E40").Value = '  +0.78%  '
$ws.Range("D41").Value = '0.2229'
$ws.Range("E41").Value = '  +5.91%  '
$ws.Range("D42").Value = '1.279'
$ws.Range("E42").Value = '  +1.89%  '
$ws.Range("D43").Value = '0.6450'
$ws.Range("E43").Value = '  +5.06%  '
$ws.Range("D44").Value = '0.9952'
$ws.Range("E44").Value = '  -0.20%  '
$ws.Range("D45").Value = '14.00'
$ws.Range("E45").Value = '  +6.38%  '
$ws.Range("D46").Value = '0.6142'
$ws.Range("E46").Value = '  +3.53%  '
$ws.Range("D47").Value = '3.815'
$ws.Range("E47").Value = '  +0.83%  '
$ws.Range("D48").Value = '2.124'
$ws.Range("E48").Value = '  +4.74%  '
$ws.Range("D49").Value = '129.64'
$ws.Range("E49").Value = '  +1.75%  '
$ws.Range("D50").Value = '0.07296'
$ws.Range("E50").Value = '  +1.44%  '
$ws.Range("D51").Value = '79.72'
$ws.Range("E51").Value = '  +5.32%  '
